$wb = $excel.ActiveWorkbook

$ALC = $wb.Worksheets.Item("ALC")
$ARM = $wb.Worksheets.Item("ARM")
$BSM = $wb.Worksheets.Item("BSM")
$CRP = $wb.Worksheets.Item("CRP")
$CUL = $wb.Worksheets.Item("CUL")
$GSM = $wb.Worksheets.Item("GSM")
$LTW = $wb.Worksheets.Item("LTW")
$WVR = $wb.Worksheets.Item("WVR")

$ALC.Range("H33").Value = 225.29167
$ALC.Range("I33").Value = 195.83333
$ALC.Range("K33").Value = 195.83333
$ALC.Range("M33").Value = 33.16667000000001
$ALC.Range("H47").Value = 5179.8
$ALC.Range("I47").Value = 5179.8
$ALC.Range("K47").Value = 5179.8
$ALC.Range("M47").Value = -4207.8
$ALC.Range("H49").Value = 206.33333
$ALC.Range("I49").Value = 207
$ALC.Range("K49").Value = 621
$ALC.Range("M49").Value = -485
$ALC.Range("H54").Value = 0
$ALC.Range("I54").Value = 0
$ALC.Range("K54").Value = 0
$ALC.Range("M54").ClearContents()
$ALC.Range("H100").Value = 4856.1055
$ALC.Range("I100").Value = 1795.2727
$ALC.Range("J100").Value = 9064.75
$ALC.Range("K100").Value = 1795.2727
$ALC.Range("L100").Value = 9064.75
$ALC.Range("M100").Value = -1254.2727
$ALC.Range("N100").Value = -10146.75
$ALC.Range("H108").Value = 78950
$ALC.Range("J108").Value = 78950
$ALC.Range("L108").Value = 78950
$ALC.Range("N108").Value = -86630
$ALC.Range("H113").Value = 3659.8
$ALC.Range("I113").Value = 4199.6665
$ALC.Range("J113").Value = 2850
$ALC.Range("K113").Value = 4199.6665
$ALC.Range("L113").Value = 2850
$ALC.Range("M113").Value = -945.6665000000003
$ALC.Range("N113").Value = -9358
$ALC.Range("H116").Value = 35824.855
$ALC.Range("I116").Value = 57688.668
$ALC.Range("J116").Value = 6673.1113
$ALC.Range("K116").Value = 57688.668
$ALC.Range("L116").Value = 6673.1113
$ALC.Range("M116").Value = -54246.668
$ALC.Range("N116").Value = -13557.1113
$ALC.Range("H125").Value = 3033.3333
$ALC.Range("J125").Value = 3487.5
$ALC.Range("L125").Value = 31387.5
$ALC.Range("N125").Value = -36307.5
$ALC.Range("H130").Value = 135292.5
$ALC.Range("J130").Value = 135292.5
$ALC.Range("L130").Value = 135292.5
$ALC.Range("N130").Value = -145332.5
$ALC.Range("H132").Value = 2269.5386
$ALC.Range("I132").Value = 2200.44
$ALC.Range("K132").Value = 6601.32
$ALC.Range("M132").Value = -4071.32
$ALC.Range("H137").Value = 2202.4075
$ALC.Range("I137").Value = 1349
$ALC.Range("J137").Value = 3121.4614
$ALC.Range("K137").Value = 4047
$ALC.Range("L137").Value = 9364.3842
$ALC.Range("M137").Value = -1497
$ALC.Range("N137").Value = -14464.3842
$ARM.Range("H32").Value = 2545.5
$ARM.Range("I32").Value = 2545.5
$ARM.Range("J32").Value = 0
$ARM.Range("K32").Value = 2545.5
$ARM.Range("L32").Value = 0
$ARM.Range("M32").ClearContents()
$ARM.Range("N32").Value = -2258.5
$ARM.Range("H61").Value = 3703.7666
$ARM.Range("I61").Value = 1805.3636
$ARM.Range("K61").Value = 1805.3636
$ARM.Range("M61").Value = -1593.3636
$ARM.Range("H97").Value = 692.80554
$ARM.Range("I97").Value = 620.4815
$ARM.Range("K97").Value = 620.4815
$ARM.Range("M97").Value = -124.4815
$ARM.Range("H98").Value = 71162.89
$ARM.Range("J98").Value = 71162.89
$ARM.Range("L98").Value = 71162.89
$ARM.Range("N98").Value = -77152.89
$ARM.Range("H122").Value = 2865.2
$ARM.Range("I122").Value = 2552.08
$ARM.Range("J122").Value = 3256.6
$ARM.Range("K122").Value = 7656.24
$ARM.Range("L122").Value = 9769.799999999999
$ARM.Range("M122").Value = -5206.24
$ARM.Range("N122").Value = -14669.8
$ARM.Range("H132").Value = 2161.4722
$ARM.Range("I132").Value = 2062.5508
$ARM.Range("K132").Value = 6187.6524
$ARM.Range("M132").Value = -3657.6524
$ARM.Range("H136").Value = 3703.7666
$ARM.Range("I136").Value = 1805.3636
$ARM.Range("K136").Value = 5416.0908
$ARM.Range("M136").Value = -2866.0908
$BSM.Range("H134").Value = 3889.8462
$BSM.Range("I134").Value = 1955.32
$BSM.Range("K134").Value = 5865.96
$BSM.Range("M134").Value = -3330.96
$CRP.Range("H58").Value = 3788
$CRP.Range("I58").Value = 1673.6666
$CRP.Range("K58").Value = 1673.6666
$CRP.Range("M58").Value = -1470.6666
$CRP.Range("H74").Value = 46299.332
$CRP.Range("I74").Value = 44585
$CRP.Range("J74").Value = 47156.5
$CRP.Range("K74").Value = 44585
$CRP.Range("L74").Value = 47156.5
$CRP.Range("M74").Value = -43711
$CRP.Range("N74").Value = -48904.5
$CRP.Range("H77").Value = 46299.332
$CRP.Range("I77").Value = 44585
$CRP.Range("J77").Value = 47156.5
$CRP.Range("K77").Value = 133755
$CRP.Range("L77").Value = 141469.5
$CRP.Range("M77").Value = -129387
$CRP.Range("N77").Value = -150205.5
$CRP.Range("H132").Value = 1774.0358
$CRP.Range("I132").Value = 1563.8077
$CRP.Range("J132").Value = 4507
$CRP.Range("K132").Value = 4691.4231
$CRP.Range("L132").Value = 13521
$CRP.Range("M132").Value = -2161.4231
$CRP.Range("N132").Value = -18581
$CRP.Range("H134").Value = 2992.8572
$CRP.Range("I134").Value = 2882.16
$CRP.Range("K134").Value = 8646.48
$CRP.Range("M134").Value = -6111.48
$CRP.Range("H136").Value = 3788
$CRP.Range("I136").Value = 1673.6666
$CRP.Range("K136").Value = 5020.9998
$CRP.Range("M136").Value = -2470.9998
$CRP.Range("H139").Value = 89997
$CRP.Range("I139").Value = 0
$CRP.Range("K139").Value = 0
$CRP.Range("M139").ClearContents()
$CRP.Range("H141").Value = 544997.25
$CRP.Range("J141").Value = 544997.25
$CRP.Range("L141").Value = 544997.25
$CRP.Range("N141").Value = -555357.25
$CUL.Range("H7").Value = 1312628.1
$CUL.Range("J7").Value = 4200015
$CUL.Range("L7").Value = 12600045
$CUL.Range("N7").Value = -12600269
$CUL.Range("H113").Value = 829.5
$CUL.Range("J113").Value = 919.6667
$CUL.Range("L113").Value = 2759.0001
$CUL.Range("N113").Value = -7099.0001
$GSM.Range("H97").Value = 3012.4736
$GSM.Range("I97").Value = 1680.5
$GSM.Range("J97").Value = 6742
$GSM.Range("K97").Value = 1680.5
$GSM.Range("L97").Value = 6742
$GSM.Range("M97").Value = -1184.5
$GSM.Range("N97").Value = -7734
$GSM.Range("H107").Value = 1819.5834
$GSM.Range("I107").Value = 1259.1111
$GSM.Range("J107").Value = 3501
$GSM.Range("K107").Value = 1259.1111
$GSM.Range("L107").Value = 3501
$GSM.Range("M107").Value = 660.8888999999999
$GSM.Range("N107").Value = -7341
$GSM.Range("H122").Value = 4400.5947
$GSM.Range("I122").Value = 2434.8635
$GSM.Range("J122").Value = 7283.6665
$GSM.Range("K122").Value = 7304.5905
$GSM.Range("L122").Value = 21850.9995
$GSM.Range("M122").Value = -4854.5905
$GSM.Range("N122").Value = -26750.9995
$GSM.Range("H123").Value = 74996.25
$GSM.Range("J123").Value = 74996.25
$GSM.Range("L123").Value = 74996.25
$GSM.Range("N123").Value = -79896.25
$GSM.Range("H126").Value = 2604.7144
$GSM.Range("I126").Value = 1968.3
$GSM.Range("J126").Value = 2958.2778
$GSM.Range("K126").Value = 5904.9
$GSM.Range("L126").Value = 8874.8334
$GSM.Range("M126").Value = -3434.9
$GSM.Range("N126").Value = -13814.8334
$GSM.Range("H132").Value = 1426.0613
$GSM.Range("I132").Value = 1325.4884
$GSM.Range("J132").Value = 2146.8333
$GSM.Range("K132").Value = 3976.4652
$GSM.Range("L132").Value = 6440.499899999999
$GSM.Range("M132").Value = -1446.4652
$GSM.Range("N132").Value = -11500.4999
$LTW.Range("H22").Value = 53949.316
$LTW.Range("I22").Value = 126254.75
$LTW.Range("J22").Value = 1363.5454
$LTW.Range("K22").Value = 126254.75
$LTW.Range("L22").Value = 1363.5454
$LTW.Range("M22").Value = -125959.75
$LTW.Range("N22").Value = -1953.5454
$LTW.Range("H27").Value = 53949.316
$LTW.Range("I27").Value = 126254.75
$LTW.Range("J27").Value = 1363.5454
$LTW.Range("K27").Value = 126254.75
$LTW.Range("L27").Value = 1363.5454
$LTW.Range("M27").Value = -126147.75
$LTW.Range("N27").Value = -1577.5454
$LTW.Range("H132").Value = 2524.453
$LTW.Range("I132").Value = 2429.587
$LTW.Range("J132").Value = 3147.8572
$LTW.Range("K132").Value = 7288.761
$LTW.Range("L132").Value = 9443.571599999999
$LTW.Range("M132").Value = -4758.761
$LTW.Range("N132").Value = -14503.5716
$WVR.Range("H122").Value = 2693.1924
$WVR.Range("I122").Value = 2575.9
$WVR.Range("J122").Value = 3084.1667
$WVR.Range("K122").Value = 7727.700000000001
$WVR.Range("L122").Value = 9252.500100000001
$WVR.Range("M122").Value = -5277.700000000001
$WVR.Range("N122").Value = -14152.5001
$WVR.Range("H123").Value = 54533.332
$WVR.Range("J123").Value = 54533.332
$WVR.Range("L123").Value = 54533.332
$WVR.Range("N123").Value = -64333.332
$WVR.Range("H136").Value = 8425766
$WVR.Range("I136").Value = 10794660
$WVR.Range("J136").Value = 3031.7778
$WVR.Range("K136").Value = 32383980
$WVR.Range("L136").Value = 9095.3334
$WVR.Range("M136").Value = -32381430
$WVR.Range("N136").Value = -14195.3334
